# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1024
$wsExhibit.Range("F3").Value = 2140
$wsExhibit.Range("F4").Value = 3
$wsExhibit.Range("F5").Value = 468

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1024
$wsAll.Range("F5").Value = 2140
$wsAll.Range("F6").Value = 3
$wsAll.Range("F7").Value = 468
